# Apply updates to Sheet1 matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("N2").Value = 19
$ws.Range("O2").Value = 1.13
$ws.Range("P2").Value = 6
$ws.Range("W2").Value = 12
$ws.Range("X2").Value = 11
$ws.Range("AF2").Value = 34
$ws.Range("BB2").Value = 67

# Row 3 updates
$ws.Range("G3").Value = 2.7
$ws.Range("L3").Value = 3.25
$ws.Range("R3").Value = 1.75
$ws.Range("U3").Value = 1.8
$ws.Range("V3").Value = 1.91
$ws.Range("AC3").Value = 9
$ws.Range("AN3").Value = 4.75
